$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A55").Value = "FlaskMegaTut"
$ws.Range("B55").Value = 43378
$ws.Range("C55").Value = 0.034722222222222224

$ws.Range("B54:C54").Copy()
$ws.Range("B55:C55").PasteSpecial(-4122)

$ws.Range("A56").Value = "Octoparse"

$ws.Range("B56").Select()
$ws.Application.ActiveWindow.ScrollRow = 31
